$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Append a new machine record (row 33), following the same column layout
# used by the existing rows: id, name, mac_address, serial_num, ip_address,
# mspec_id, lang_code, is_active, cr_by, cr_dtimes, eff_dtimes
$newRow = 33

$ws.Cells.Item($newRow, 1).Value = 10032
$ws.Cells.Item($newRow, 2).Value = "Machine 32"
$ws.Cells.Item($newRow, 3).Value = "F4-30-B9-D4-CD-6F"
$ws.Cells.Item($newRow, 4).Value = "FB5962911665"
$ws.Cells.Item($newRow, 5).Value = "192.168.0.358"
$ws.Cells.Item($newRow, 6).Value = 1001
$ws.Cells.Item($newRow, 7).Value = "eng"
$ws.Cells.Item($newRow, 8).Value = $true
$ws.Cells.Item($newRow, 9).Value = "superadmin"
$ws.Cells.Item($newRow, 10).Value = "now()"
$ws.Cells.Item($newRow, 11).Value = "now()"

# Restore the view/selection state saved with the workbook: scrolled so
# row 22 is at the top, with C28 as the active (selected) cell.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C28").Select()
